$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column F ("Voltage") to hold the new
# "Dielectric" column. This shifts the existing F:N columns to G:O.
$ws.Columns("F").EntireColumn.Insert()

# Match the column width used by the author for the new column.
$ws.Columns("F").ColumnWidth = 8.6

# Header
$ws.Range("F1").Value = "Dielectric"

# Dielectric values per row, matching the Description (column B) for each part.
$ws.Range("F2").Value = "C0G"
$ws.Range("F3").Value = "C0G"
$ws.Range("F4").Value = "C0G"
$ws.Range("F5").Value = "C0G"
$ws.Range("F6").Value = "C0G"
$ws.Range("F7").Value = "C0G"
$ws.Range("F8").Value = "C0G"
$ws.Range("F9").Value = "C0G"
$ws.Range("F10").Value = "C0G"
$ws.Range("F11").Value = "X7R"
$ws.Range("F12").Value = "X7R"
$ws.Range("F13").Value = "X7R"
$ws.Range("F14").Value = "X7R"
$ws.Range("F15").Value = "X7R"
$ws.Range("F16").Value = "X7R"
$ws.Range("F17").Value = "X7R"
$ws.Range("F18").Value = "X7R"
$ws.Range("F19").Value = "X5R"
$ws.Range("F20").Value = "X5R"
$ws.Range("F21").Value = "X5R"
$ws.Range("F22").Value = "X5R"
